# Add delay to export
# Insert a header row (Patient / Emeliza Yabut) at the top of the
# "Health Habits" sheet, matching the row already present on "Daily Trackers".
# Row 1 on this sheet is currently empty (data starts at row 2), so we just
# populate the cells directly rather than shifting existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Health Habits")

$ws.Cells.Item(1, 1).Value = "Patient"
$ws.Cells.Item(1, 1).Font.Bold = $true

$ws.Cells.Item(1, 2).Value = "Emeliza Yabut"
$ws.Cells.Item(1, 2).Font.Bold = $false
